$d = $word.ActiveDocument

# Locate the run that reads "Lines 6-16<TAB>Install required packages" and
# determine the exact character range covering the <tab> + the text that
# follows it ("Install required packages").
$labelRng = $d.Content.Duplicate
$labelRng.Find.Execute("Lines 6-16", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

$textRng = $d.Content.Duplicate
$textRng.Find.Execute("Install required packages", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)

$targetRng = $d.Range($labelRng.End, $textRng.End)

# Replace that run's contents with three runs: a tab-only run, a run
# containing "Load", and a run containing " required packages" - turning
# "Install required packages" into "Load required packages" while keeping
# the tab in its own run, matching how Word split the edit.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:tab/></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t>Load</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> required packages</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$targetRng.InsertXML($xml)
